$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; this shifts existing rows 32:45 down to 33:46
# and preserves their values/styles (including column D's date style).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value = 4
$ws.Cells.Item(32, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(32, 3).Value = "Los Lagos"
$ws.Cells.Item(32, 4).Value = 44489
$ws.Cells.Item(32, 5).Value = 10
$ws.Cells.Item(32, 6).Value = 100112026
$ws.Cells.Item(32, 7).Value = "Haba"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 35
$ws.Cells.Item(32, 11).Value = 11000
$ws.Cells.Item(32, 12).Value = 11000
$ws.Cells.Item(32, 13).Value = 11000
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región Metropolitana"
$ws.Cells.Item(32, 16).Value = 440
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
